$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "flag_do" section (and everything below it) down by one row
# to make room for a new "sensor replaced (prefaced by depth of sensor)"
# row right under the flag_temp section.
$ws.Rows("9:9").Insert()

# Insert three new rows for the additional flag_do codes, pushing
# flag_met and everything below it down.
$ws.Rows("13:15").Insert()

# New row for the flag_do section: presumed calibrated
$ws.Range("B13").Value = "presumed calibrated"
$ws.Range("A13").Value = "pc"

# New row under flag_temp: sensor replaced (prefaced by depth of sensor)
$ws.Range("A8").Value = "r"
$ws.Range("B8").Value = "sensor replaced (prefaced by depth of sensor)"

# New row for the flag_do section: sensor replaced
$ws.Range("A14").Value = "r"
$ws.Range("B14").Value = "sensor replaced"

# New row for the flag_do section: sensor failed
$ws.Range("A15").Value = "f"
$ws.Range("B15").Value = "sensor failed"

# New row for the flag_met section: precip values interpolated
$ws.Range("A21").Value = "i"
$ws.Range("B21").Value = "precip values interpolated (linear) from adjacent records"

$ws.Range("B20").Select()
